$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row (row 1) cells to the new snake_case header names.
$ws.Cells.Item(1, 1).Value = "model"
$ws.Cells.Item(1, 2).Value = "release_time"
$ws.Cells.Item(1, 3).Value = "CPU_model"
$ws.Cells.Item(1, 4).Value = "CPU_frequency"
$ws.Cells.Item(1, 5).Value = "memory"
$ws.Cells.Item(1, 6).Value = "storage"
$ws.Cells.Item(1, 7).Value = "version"
$ws.Cells.Item(1, 8).Value = "user"
$ws.Cells.Item(1, 9).Value = "prevalence_bopt"
$ws.Cells.Item(1, 10).Value = "frequency_bopt"
$ws.Cells.Item(1, 11).Value = "rate_bopt"
$ws.Cells.Item(1, 12).Value = "energy_bopt"
$ws.Cells.Item(1, 13).Value = "prevalence_aopt"
$ws.Cells.Item(1, 14).Value = "frequency_aopt"
$ws.Cells.Item(1, 15).Value = "rate_aopt"
$ws.Cells.Item(1, 16).Value = "energy_aopt"

# Widen columns I and J to match the new header text lengths.
$ws.Columns.Item(9).ColumnWidth = 34.6666666666667
$ws.Columns.Item(10).ColumnWidth = 50.5

# Move the active selection to R14 (was U80).
$ws.Range("R14").Select()
